$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A6").Value = "null_threshhold"
$ws.Range("B6").Value = 1
